$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Dear ," -> "Dear sean,"  (fill in the empty bold run between
#    "Dear " and "," with the customer's name, keeping the 3-run
#    structure: "Dear " / bold "sean" / ",")
# ------------------------------------------------------------------
$pDear = $d.Paragraphs(6)
$rDear = $pDear.Range
$dearXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Dear </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>sean</w:t></w:r><w:r><w:t>,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rDear.InsertXML($dearXml)

# ------------------------------------------------------------------
# 2. Purchase date: 2022/12/04 -> 2022/12/11
# ------------------------------------------------------------------
$pDate = $d.Paragraphs(7)
$rDate = $d.Range($pDate.Range.Start, $pDate.Range.End)
$rDate.Find.Execute("2022/12/04", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rDate.Text = "2022/12/11"

# ------------------------------------------------------------------
# 3. Invoice number: 1 -> 25
# ------------------------------------------------------------------
$pInv = $d.Paragraphs(8)
$rInv = $d.Range($pInv.Range.Start, $pInv.Range.End)
$rInv.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rInv.Text = "25"

# ------------------------------------------------------------------
# 4. Remove one of the two blank paragraphs that follow the table
#    (there were two <w:p/> there, now there should be only one).
# ------------------------------------------------------------------
$d.Paragraphs(33).Range.Delete()

# ------------------------------------------------------------------
# 5. Total price: 150.82 -> 43.09
# ------------------------------------------------------------------
$pTotal = $d.Paragraphs(34)
$rTotal = $d.Range($pTotal.Range.Start, $pTotal.Range.End)
$rTotal.Find.Execute("150.82", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rTotal.Text = "43.09"

# ------------------------------------------------------------------
# 6. Table: drop the "Fallout 4" line item, replace the "LEGO Star
#    Wars" line item with "Mario Kart Wii" and renumber/repoint it
#    into what was previously the first data row, then delete the
#    now-redundant second data row. (Table operations are done last
#    because touching Table/Cell objects invalidates later
#    Paragraphs(...) index lookups in this host.)
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Cell(2, 1).Range.Text = "2"
$t.Cell(2, 2).Range.Text = " Mario Kart Wii  - Wii"
$t.Cell(2, 3).Range.Text = $t.Cell(3, 3).Range.Text
$t.Cell(2, 4).Range.Text = $t.Cell(3, 4).Range.Text
$t.Cell(2, 5).Range.Text = $t.Cell(3, 5).Range.Text
$t.Cell(2, 6).Range.Text = $t.Cell(3, 6).Range.Text
$t.Rows(3).Delete()
